# Weekly update: a new record is inserted at the top of the data table
# (row 118), pushing all subsequent records down by one row and
# appending the previously-last record (old row 140) as the new row 141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 118, shifting rows 118:140 down
# to 119:141. This grows the used range from A1:R140 to A1:R141.
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with this week's record.
$ws.Cells.Item(118, 1).Value = 8
$ws.Cells.Item(118, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 45275
$ws.Cells.Item(118, 5).Value = 4
$ws.Cells.Item(118, 6).Value = 100112030
$ws.Cells.Item(118, 7).Value = "Poroto granado"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 500
$ws.Cells.Item(118, 11).Value = 39000
$ws.Cells.Item(118, 12).Value = 40000
$ws.Cells.Item(118, 13).Value = 39500
$ws.Cells.Item(118, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(118, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(118, 16).Value = 1580
$ws.Cells.Item(118, 17).Value = 25
$ws.Cells.Item(118, 18).Value = "Hortaliza"
